$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Joseph Mensah: update email
$ws.Range("G2").Value = "jmesnah@gmail.com"

# Row 3 - Njoroge Kariuki: update email and country
$ws.Range("G3").Value = "kariukij@gmail.com"
$ws.Range("J3").Value = "Niger"

# Row 4 - Kofinuel Emmanuel: update email and country
$ws.Range("G4").Value = "emamnuel@gmail.com"
$ws.Range("J4").Value = "Niger"

# Row 5 - Mary Jane Mensah: update email and thematic field
$ws.Range("G5").Value = "janemama@gmail.com"
$ws.Range("K5").Value = "Education"

# Row 6 - Enyonam Mensah: update email and institution
$ws.Range("G6").Value = "etonammensah@gmail.com"
$ws.Range("H6").Value = "Université Abdou Moumouni`t"

# Row 7 - George KuntaKunta: update email and ACE
$ws.Range("G7").Value = "georgiekuntae@gmail.com"
$ws.Range("I7").Value = "MS4SSA"

# Update the active cell selection to L11
[void]$ws.Range("L11").Select()
